$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.202.39'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.906.45'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.89%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.01'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.09%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5254'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3785'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07279'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.31%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.65%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07689'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.914.46'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '95.07'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.277'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.01%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008625'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.39%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.45%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.266.61'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.074'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.144.44'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.64'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.442'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.315'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +11.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '145.82'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.81%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.15'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.57%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.736'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.81'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.967'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.817'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09223'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.88%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8163'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +9.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05067'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.242'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +7.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.995'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.310'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.594'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5688'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.64%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01992'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.076'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.990'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +5.06%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.632'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '119.29'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1517'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4839'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.19%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.626'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +4.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.60'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.73'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.05%  '
